# Update the "timestamp" column (Z) for every data row (2-112) in the
# Log_Muestras sheet, replacing the old run timestamps with the new ones
# recorded for this (re-)execution of the notebook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of contiguous row-ranges -> new timestamp string (grouped from the
# per-row values so we can set them with as few COM calls as possible).
$groups = @(
    @{ Start = 2;   End = 7;   Value = "2025-10-17T07:09:36.653913" },
    @{ Start = 8;   End = 17;  Value = "2025-10-17T07:09:36.654913" },
    @{ Start = 18;  End = 28;  Value = "2025-10-17T07:09:36.655916" },
    @{ Start = 29;  End = 34;  Value = "2025-10-17T07:09:36.656913" },
    @{ Start = 35;  End = 44;  Value = "2025-10-17T07:09:36.658124" },
    @{ Start = 45;  End = 45;  Value = "2025-10-17T07:09:36.659123" },
    @{ Start = 46;  End = 48;  Value = "2025-10-17T07:09:36.711797" },
    @{ Start = 49;  End = 54;  Value = "2025-10-17T07:09:36.712796" },
    @{ Start = 55;  End = 60;  Value = "2025-10-17T07:09:36.713796" },
    @{ Start = 61;  End = 66;  Value = "2025-10-17T07:09:36.714795" },
    @{ Start = 67;  End = 72;  Value = "2025-10-17T07:09:36.715796" },
    @{ Start = 73;  End = 74;  Value = "2025-10-17T07:09:36.716795" },
    @{ Start = 75;  End = 75;  Value = "2025-10-17T07:09:36.771438" },
    @{ Start = 76;  End = 85;  Value = "2025-10-17T07:09:36.772438" },
    @{ Start = 86;  End = 90;  Value = "2025-10-17T07:09:36.773434" },
    @{ Start = 91;  End = 94;  Value = "2025-10-17T07:09:36.774437" },
    @{ Start = 95;  End = 100; Value = "2025-10-17T07:09:36.775438" },
    @{ Start = 101; End = 102; Value = "2025-10-17T07:09:36.776436" },
    @{ Start = 103; End = 112; Value = "2025-10-17T07:09:36.829443" }
)

foreach ($g in $groups) {
    for ($r = $g.Start; $r -le $g.End; $r++) {
        $ws.Range("Z$r").Value = $g.Value
    }
}
